$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update tornado-plot example data values
$ws.Range("C2").Value = -0.5
$ws.Range("C3").Value = -0.2

$ws.Range("B8").Value = -0.5

$ws.Range("B9").Value = 0.1

$ws.Range("B11").Value = -4
$ws.Range("C11").Value = -2
$ws.Range("D11").Value = -1.5

$ws.Range("B14").Value = 2

$ws.Range("B17").Value = -2
$ws.Range("C17").Value = -0.4
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 3

$ws.Range("C18").Value = -0.7
$ws.Range("G18").Value = 5

# Update selection and zoom to match the final workbook view state
$ws.Range("E28").Select()
$excel.ActiveWindow.Zoom = 100
